# Roll the yearly income-statement table forward by one fiscal year:
# drop the oldest column (1396/12), shift the remaining four columns
# left (D<-E, E<-F, F<-G, G<-H), and populate the newly freed last
# column (H) with the latest period's figures / publish date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-12-18 (3)"
$ws.Range("E9").Value = "1400-04-14 (8)"
$ws.Range("F9").Value = "1401-04-12 (11)"
$ws.Range("G9").Value = "1401-10-28 (7)"
$ws.Range("H9").Value = "1402-02-25 (8)"

$ws.Range("D11").Value = 29687
$ws.Range("E11").Value = 51198
$ws.Range("F11").Value = 55842
$ws.Range("G11").Value = 65226
$ws.Range("H11").Value = 77250

$ws.Range("D12").Value = -16490
$ws.Range("E12").Value = -25510
$ws.Range("F12").Value = -28606
$ws.Range("G12").Value = -34901
$ws.Range("H12").Value = -42971

$ws.Range("D13").Value = 13196
$ws.Range("E13").Value = 25688
$ws.Range("F13").Value = 27236
$ws.Range("G13").Value = 30325
$ws.Range("H13").Value = 34279

$ws.Range("D14").Value = -7909
$ws.Range("E14").Value = -8549
$ws.Range("F14").Value = -7918
$ws.Range("G14").Value = -11691
$ws.Range("H14").Value = -13504

$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"
$ws.Range("G16").Value = "-"
$ws.Range("H16").Value = "-"

$ws.Range("D17").Value = 5291
$ws.Range("E17").Value = 17139
$ws.Range("F17").Value = 19318
$ws.Range("G17").Value = 18634
$ws.Range("H17").Value = 20775

$ws.Range("D18").Value = -3184
$ws.Range("E18").Value = -3227
$ws.Range("F18").Value = -2286
$ws.Range("G18").Value = -2555
$ws.Range("H18").Value = -2943

$ws.Range("D19").Value = -984
$ws.Range("E19").Value = -281
$ws.Range("F19").Value = -220
$ws.Range("G19").Value = 370
$ws.Range("H19").Value = 287

$ws.Range("D20").Value = 1122
$ws.Range("E20").Value = 13631
$ws.Range("F20").Value = 16812
$ws.Range("G20").Value = 16449
$ws.Range("H20").Value = 18119

$ws.Range("D21").Value = -242
$ws.Range("E21").Value = -1928
$ws.Range("F21").Value = -849
$ws.Range("G21").Value = -826
$ws.Range("H21").Value = -725

$ws.Range("D22").Value = 880
$ws.Range("E22").Value = 11702
$ws.Range("F22").Value = 15962
$ws.Range("G22").Value = 15623
$ws.Range("H22").Value = 17394

$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

$ws.Range("D24").Value = 880
$ws.Range("E24").Value = 11702
$ws.Range("F24").Value = 15962
$ws.Range("G24").Value = 15623
$ws.Range("H24").Value = 17394

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("D26").Value = 4943
$ws.Range("E26").Value = 3898
$ws.Range("F26").Value = 6634
$ws.Range("G26").Value = 5685
$ws.Range("H26").Value = 11335

$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
